# Applies crypto price/volume updates from the GitHub Actions refresh
# (commit: "Updated cryptos list on Mon Nov 11 18:51:37 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "85.321.56"
$ws.Range("E2").Value = "  +5.64%  "

$ws.Range("D3").Value = "3.319.43"
$ws.Range("E3").Value = "  +2.47%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.93"
$ws.Range("E5").Value = "  +2.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "637.16"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.323"
$ws.Range("E7").Value = "  +10.24%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  -2.39%  "

$ws.Range("D10").Value = "3.322.79"
$ws.Range("E10").Value = "  +2.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.598"
$ws.Range("E11").Value = "  -2.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000276"
$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").Value = "3.916.86"
$ws.Range("E14").Value = "  +2.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.30"
$ws.Range("E15").Value = "  +3.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  -1.41%  "

$ws.Range("D17").Value = "85.001.96"
$ws.Range("E17").Value = "  +5.45%  "

$ws.Range("D18").Value = "3.315.11"
$ws.Range("E18").Value = "  +2.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.61"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.17"
$ws.Range("E20").Value = "  +3.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "440.46"
$ws.Range("E21").Value = "  -2.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.18"
$ws.Range("E22").Value = "  -2.77%  "

$ws.Range("E23").Value = "  -3.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.40"
$ws.Range("E24").Value = "  +3.70%  "

$ws.Range("E25").Value = "  +12.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.22"
$ws.Range("E26").Value = "  +9.94%  "

$ws.Range("D27").Value = "3.473.17"
$ws.Range("E27").Value = "  +2.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "78.17"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000131"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "606.73"
$ws.Range("E31").Value = "  +6.40%  "

$ws.Range("E32").Value = "  +30.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.23"
$ws.Range("E33").Value = "  -1.92%  "

$ws.Range("E34").Value = "  +0.35%  "

$ws.Range("E35").Value = "  +1.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.05"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.152"
$ws.Range("E37").Value = "  -4.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.21"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.41"
$ws.Range("E39").Value = "  +8.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.419"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.22"
$ws.Range("E42").Value = "  +4.27%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("E43").Value = "  +8.31%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.03"
$ws.Range("E44").Value = "  +9.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "160.00"
$ws.Range("E45").Value = "  -3.08%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "189.71"
$ws.Range("E47").Value = "  -2.23%  "

$ws.Range("E48").Value = "  -0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.89"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.789"
$ws.Range("E50").Value = "  -2.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.51"
$ws.Range("E51").Value = "  +1.89%  "

